$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 header: was a numeric 0, now should be the text "Description"
$ws.Range("A1").Value = "Description"

# Update the "Move to location (x, y) ..." sentences in column A (rows 2-11)
# with new coordinates, keeping the rest of each sentence unchanged.
$ws.Range("A2").Value  = "Move to location (11, 8) and remove the toolkit."
$ws.Range("A3").Value  = "Move to location (7, 5) and remove the liquid spill."
$ws.Range("A4").Value  = "Move to location (8, 6) and remove the large debris."
$ws.Range("A5").Value  = "Move to location (2, 4) and remove the dust."
$ws.Range("A6").Value  = "Move to location (5, 2) and remove the grass."
$ws.Range("A7").Value  = "Move to location (6, 7) and remove the small debris."
$ws.Range("A8").Value  = "Move to location (3, 6) and remove the vehicle."
$ws.Range("A9").Value  = "Move to location (6, 6) and remove the construction materials."
$ws.Range("A10").Value = "Move to location (3, 9) and remove the tree branches."
$ws.Range("A11").Value = "Move to location (6, 6) and remove the screws."
